$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 4: student id changed, name/surname stay the same
$ws.Range("A4").Value = 181004

# Add two new students in rows 6 and 7
$ws.Range("A6").Value = 186037
$ws.Range("B6").Value = "Дарко"
$ws.Range("C6").Value = "Ристевски"

$ws.Range("A7").Value = 183160
$ws.Range("B7").Value = "Јована"
$ws.Range("C7").Value = "Ѓурковска"

# Move selection to A5, matching the saved view state in the diff
$ws.Range("A5").Select()
